$wb = $excel.ActiveWorkbook

# --- Sheet "About": update the two explanatory notes from 100% to 90% ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A16").Value = "For industries, we assign 90% to all industries except mining, agriculture, and waste management, as the activities"
$wsAbout.Range("A14").Value = "For simplicity, for the electricity sector, we estimate a share of 90% for all fossil power plant types, biomass, and MSW."

# --- Sheet "CPPbES": CO2 capture potential by electricity source, 1 -> 0.9 for applicable rows ---
$wsES = $wb.Worksheets.Item("CPPbES")
$wsES.Range("B2").Value = 0.9   # hard coal
$wsES.Range("B3").Value = 0.9   # natural gas nonpeaker
$wsES.Range("B9").Value = 0.9   # biomass
$wsES.Range("B11").Value = 0.9  # petroleum
$wsES.Range("B12").Value = 0.9  # natural gas peaker
$wsES.Range("B13").Value = 0.9  # lignite
$wsES.Range("B15").Value = 0.9  # crude oil
$wsES.Range("B16").Value = 0.9  # heavy or residual fuel oil
$wsES.Range("B17").Value = 0.9  # municipal solid waste

# --- Sheet "CPPbI": CO2 capture potential by industry, 1 -> 0.9 for applicable rows ---
$wsI = $wb.Worksheets.Item("CPPbI")
$wsI.Range("B2").Value = 0.9  # cement and other carbonates
$wsI.Range("B3").Value = 0.9  # natural gas and petroleum systems
$wsI.Range("B4").Value = 0.9  # iron and steel
$wsI.Range("B5").Value = 0.9  # chemicals
$wsI.Range("B9").Value = 0.9  # other industries
